$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"


$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C16").Value = 8

$ws.Range("D16").Value = 2

$ws.Range("E16").Value = 300

$ws.Range("F16").Value = 17

$ws.Range("G16").Value = 8

$ws.Range("H16").Value = 112.5

$ws.Range("I16").Value = 209

$ws.Range("J16").Value = 137

$ws.Range("K16").Value = 52.554744525547

$ws.Range("L16").Value = 71.311475409836

$ws.Range("M16").Value = 35.714285714285

$ws.Range("N16").Value = -81.389136242208

$ws.Range("C17").Value = 6

$ws.Range("D17").Value = 2

$ws.Range("E17").Value = 200

$ws.Range("G17").Value = 20

$ws.Range("H17").Value = 10

$ws.Range("I17").Value = 291

$ws.Range("J17").Value = 229

$ws.Range("K17").Value = 27.074235807860

$ws.Range("L17").Value = 84.177215189873

$ws.Range("M17").Value = 206.315789473684

$ws.Range("N17").Value = -23.821989528795

$ws.Range("C18").Value = 8

$ws.Range("D18").Value = 1

$ws.Range("E18").Value = 700

$ws.Range("F18").Value = 26

$ws.Range("H18").Value = 271.428571428571

$ws.Range("I18").Value = 204

$ws.Range("J18").Value = 166

$ws.Range("K18").Value = 22.891566265060

$ws.Range("L18").Value = 23.636363636363

$ws.Range("M18").Value = 98.058252427184

$ws.Range("N18").Value = -71.348314606741

$ws.Range("C19").Value = 22

$ws.Range("D19").Value = 16

$ws.Range("E19").Value = 37.5

$ws.Range("F19").Value = 67

$ws.Range("G19").Value = 49

$ws.Range("H19").Value = 36.734693877551

$ws.Range("I19").Value = 678

$ws.Range("J19").Value = 606

$ws.Range("K19").Value = 11.881188118811

$ws.Range("L19").Value = 22.162162162162

$ws.Range("M19").Value = 55.504587155963

$ws.Range("N19").Value = -29.301355578727

$ws.Range("C20").Value = 1

$ws.Range("D20").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = -50
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("F20").Value = 11

$ws.Range("G20").Value = 6

$ws.Range("H20").Value = 83.333333333333

$ws.Range("I20").Value = 87

$ws.Range("J20").Value = 67

$ws.Range("K20").Value = 29.850746268656

$ws.Range("L20").Value = 58.181818181818

$ws.Range("M20").Value = 74

$ws.Range("N20").Value = -86.532507739938

$ws.Range("C21").Value = 45

$ws.Range("D21").Value = 23

$ws.Range("E21").Value = 95.652173913043

$ws.Range("F21").Value = 144

$ws.Range("G21").Value = 93

$ws.Range("H21").Value = 54.838709677419

$ws.Range("I21").Value = 1479

$ws.Range("J21").Value = 1211

$ws.Range("K21").Value = 22.130470685384

$ws.Range("L21").Value = 38.742964352720

$ws.Range("M21").Value = 75.653206650831

$ws.Range("N21").Value = -61.594391067255

$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F22").Value = 4

$ws.Range("G22").Value = 3

$ws.Range("H22").Value = 33.333333333333

$ws.Range("I22").Value = 42

$ws.Range("K22").Value = -17.647058823529

$ws.Range("L22").Value = -25

$ws.Range("M22").Value = -32.258064516129

$ws.Range("D23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("F23").Value = 3

$ws.Range("G23").Value = 2

$ws.Range("H23").Value = 50

$ws.Range("I23").Value = 43

$ws.Range("K23").Value = 2.380952380952

$ws.Range("L23").Value = 30.303030303030

$ws.Range("M23").Value = 79.166666666666

$ws.Range("C24").Value = 51

$ws.Range("D24").Value = 21

$ws.Range("E24").Value = 142.857142857143

$ws.Range("F24").Value = 193

$ws.Range("G24").Value = 132

$ws.Range("H24").Value = 46.212121212121

$ws.Range("I24").Value = 1937

$ws.Range("J24").Value = 1779

$ws.Range("K24").Value = 8.881394041596

$ws.Range("L24").Value = 14.479905437352

$ws.Range("M24").Value = 47.188449848024

$ws.Range("C25").Value = 43

$ws.Range("D25").Value = 25

$ws.Range("E25").Value = 72

$ws.Range("F25").Value = 158

$ws.Range("H25").Value = 25.396825396825

$ws.Range("I25").Value = 1719

$ws.Range("J25").Value = 1592

$ws.Range("K25").Value = 7.977386934673

$ws.Range("L25").Value = 22.785714285714

$ws.Range("C26").Value = 13

$ws.Range("D26").Value = 8

$ws.Range("E26").Value = 62.5

$ws.Range("G26").Value = 41

$ws.Range("H26").Value = 24.390243902439

$ws.Range("I26").Value = 423

$ws.Range("J26").Value = 363

$ws.Range("K26").Value = 16.528925619834

$ws.Range("L26").Value = 31.775700934579

$ws.Range("M26").Value = 34.713375796178

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F28").Value = 3

$ws.Range("G28").Value = 4

$ws.Range("H28").Value = -25

$ws.Range("J28").Value = 53

$ws.Range("K28").Value = -3.773584905660

$ws.Range("L28").Value = -13.559322033898

$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("D31").Value = 2

$ws.Range("F31").Value = 3

$ws.Range("G31").Value = 5

$ws.Range("H31").Value = -40

$ws.Range("J31").Value = 18

$ws.Range("K31").Value = 0

$ws.Range("C33").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$ws.Range("F33").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F33").PasteSpecial(-4122)

$ws.Range("I33").Value = 1
$ws.Range("I14").Copy()
$ws.Range("I33").PasteSpecial(-4122)

$ws.Range("K33").Value = -50

$ws.Range("L33").Value = 0
